$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": add a new column BG with the "11-aug" prices
# ---------------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting of the previous header cell (BF1) onto the new header
# cell (BG1) before writing its value so it keeps the bold / bordered style.
$wsSpot.Range("BF1").Copy() | Out-Null
$wsSpot.Range("BG1").PasteSpecial(-4122) | Out-Null
$wsSpot.Range("BG1").Value = "11-aug"

$bgValues = @(34.57, 31.08, 41.46, 36.32, 38.25, 40.37, 39.06, 34.4, 50.14, 58.34, 45, 2.54, 0, 0, 0, 9, 44.5, 57.97, 90, 118.68, 135.09, 138.94, 108.87, 96)

for ($i = 0; $i -lt $bgValues.Length; $i++) {
    $row = $i + 2
    $wsSpot.Cells.Item($row, 59).Value = $bgValues[$i]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append rows 56 and 57
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

$wsGaz.Range("A56").NumberFormat = "@"
$wsGaz.Range("A56").Value = "2025-08-09"
$wsGaz.Range("A56").ClearFormats()
$wsGaz.Range("B56").Value = 31.375

$wsGaz.Range("A57").NumberFormat = "@"
$wsGaz.Range("A57").Value = "2025-08-10"
$wsGaz.Range("A57").ClearFormats()
$wsGaz.Range("B57").Value = 31.375

# ---------------------------------------------------------------------------
# Sheet "CO2": append rows 56 and 57
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A56").NumberFormat = "@"
$wsCo2.Range("A56").Value = "2025-08-09"
$wsCo2.Range("A56").ClearFormats()
$wsCo2.Range("B56").Value = 71.75

$wsCo2.Range("A57").NumberFormat = "@"
$wsCo2.Range("A57").Value = "2025-08-10"
$wsCo2.Range("A57").ClearFormats()
$wsCo2.Range("B57").Value = 71.75
